$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set purpose (column E) to "fullRNASEQ" and libraryPreparer (column B) to "H.BROWN"
# for every data row (rows 2 through 45), replacing the old "Retrofitted_2651" placeholder.
$ws.Range("E2:E45").Value = "fullRNASEQ"
$ws.Range("B2:B45").Value = "H.BROWN"

# Select B3:B45 with active cell B3, matching the saved selection state.
$ws.Range("B3:B45").Select()
